# Update forecast-error statistics (ME, MAE, MSE, RMSE, SE) for horizons Q0..Q8
# with corrected values (bug fix in upstream calculation).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Q0)
$ws.Range("B2").Value = 0.04642948062583634
$ws.Range("C2").Value = 0.2965858768643577
$ws.Range("D2").Value = 0.1211486527137855
$ws.Range("E2").Value = 0.3480641502852391
$ws.Range("F2").Value = 0.3579752231001048

# Row 3 (Q1)
$ws.Range("B3").Value = 0.07710192660348104
$ws.Range("C3").Value = 0.3275181423519703
$ws.Range("D3").Value = 0.1806004724994779
$ws.Range("E3").Value = 0.4249711431373639
$ws.Range("F3").Value = 0.4349832324714388

# Row 4 (Q2)
$ws.Range("B4").Value = 0.06899662371576315
$ws.Range("C4").Value = 0.2620176638063245
$ws.Range("D4").Value = 0.09766707149571256
$ws.Range("E4").Value = 0.3125173139135055
$ws.Range("F4").Value = 0.3183592094900543

# Row 5 (Q3)
$ws.Range("B5").Value = 0.1035751530193068
$ws.Range("C5").Value = 0.3321382172543214
$ws.Range("D5").Value = 0.1880858394500914
$ws.Range("E5").Value = 0.4336886434414572
$ws.Range("F5").Value = 0.4416942719119532

# Row 6 (Q4)
$ws.Range("B6").Value = 0.1569404967837869
$ws.Range("C6").Value = 0.3655869259556198
$ws.Range("D6").Value = 0.2335277915346266
$ws.Range("E6").Value = 0.4832471329812797
$ws.Range("F6").Value = 0.481776195163832

# Row 7 (Q5)
$ws.Range("B7").Value = 0.1649274464288797
$ws.Range("C7").Value = 0.4479266052260063
$ws.Range("D7").Value = 0.276931746074847
$ws.Range("E7").Value = 0.5262430484812574
$ws.Range("F7").Value = 0.5300443556207896

# Row 8 (Q6)
$ws.Range("B8").Value = 0.2184461953069695
$ws.Range("C8").Value = 0.4442913753428973
$ws.Range("D8").Value = 0.3428095254168975
$ws.Range("E8").Value = 0.5854993812267418
$ws.Range("F8").Value = 0.5950705354891705

# Row 9 (Q7)
$ws.Range("B9").Value = 0.04181093721508944
$ws.Range("C9").Value = 0.2983623655071265
$ws.Range("D9").Value = 0.09486659670334104
$ws.Range("E9").Value = 0.308004215398655
$ws.Range("F9").Value = 0.3737347499882843

# Row 10 (Q8)
$ws.Range("B10").Value = -0.2839541438535775
$ws.Range("C10").Value = 0.2839541438535775
$ws.Range("D10").Value = 0.08062995581161821
$ws.Range("E10").Value = 0.2839541438535775
